$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability matrix values (per commit: "added more games,
# sped up simulate game logic, and drafted optimization logic").
# Each row of the matrix represents transition probabilities from a given
# Starting_State, recomputed after simulating additional games.

# Row 2
$ws.Range("B2").Value = 0.2172523961661342
$ws.Range("C2").Value = 0.5111821086261981
$ws.Range("J2").Value = 0.01916932907348243
$ws.Range("P2").Value = 0.1597444089456869
$ws.Range("S2").Value = 0.0926517571884984
# Row 3
$ws.Range("B3").Value = 0.0245398773006135
$ws.Range("C3").Value = 0.0245398773006135
$ws.Range("J3").Value = 0.01840490797546012
$ws.Range("P3").Value = 0.7361963190184049
$ws.Range("S3").Value = 0.196319018404908
# Row 4
$ws.Range("J4").Value = 0.04347826086956522
$ws.Range("P4").Value = 0.6304347826086957
$ws.Range("S4").Value = 0.3260869565217391
# Row 6
$ws.Range("B6").Value = 0.08947368421052632
$ws.Range("F6").Value = 0.04736842105263158
$ws.Range("J6").Value = 0.2789473684210526
$ws.Range("O6").Value = 0.03157894736842105
$ws.Range("Q6").Value = 0.1736842105263158
$ws.Range("R6").Value = 0.09473684210526316
$ws.Range("S6").Value = 0.2842105263157895
# Row 7
$ws.Range("B7").Value = 0.1176470588235294
$ws.Range("D7").Value = 0.02205882352941177
$ws.Range("F7").Value = 0.06617647058823529
$ws.Range("J7").Value = 0.1838235294117647
$ws.Range("O7").Value = 0.04411764705882353
$ws.Range("Q7").Value = 0.1397058823529412
$ws.Range("R7").Value = 0.09558823529411764
$ws.Range("S7").Value = 0.3308823529411765
# Row 8
$ws.Range("B8").Value = 0.0970873786407767
$ws.Range("D8").Value = 0.01699029126213592
$ws.Range("F8").Value = 0.04368932038834952
$ws.Range("J8").Value = 0.09951456310679611
$ws.Range("O8").Value = 0.02912621359223301
$ws.Range("Q8").Value = 0.2111650485436893
$ws.Range("R8").Value = 0.1116504854368932
$ws.Range("S8").Value = 0.3907766990291262
# Row 9
$ws.Range("B9").Value = 0.06666666666666667
$ws.Range("D9").Value = 0.0380952380952381
$ws.Range("F9").Value = 0.09523809523809523
$ws.Range("J9").Value = 0.1238095238095238
$ws.Range("O9").Value = 0.009523809523809525
$ws.Range("Q9").Value = 0.2380952380952381
$ws.Range("R9").Value = 0.08095238095238096
$ws.Range("S9").Value = 0.3476190476190476
# Row 10
$ws.Range("B10").Value = 0.128337639965547
$ws.Range("D10").Value = 0.02583979328165375
$ws.Range("E10").Value = 0.003445305770887166
$ws.Range("F10").Value = 0.06546080964685616
$ws.Range("J10").Value = 0.1042204995693368
$ws.Range("O10").Value = 0.01550387596899225
$ws.Range("Q10").Value = 0.2213608957795004
$ws.Range("R10").Value = 0.0913006029285099
$ws.Range("S10").Value = 0.3445305770887166
# Row 11
$ws.Range("G11").Value = 0.1155778894472362
$ws.Range("J11").Value = 0.1407035175879397
$ws.Range("K11").Value = 0.2060301507537688
$ws.Range("L11").Value = 0.5276381909547738
$ws.Range("S11").Value = 0.01005025125628141
# Row 12
$ws.Range("G12").Value = 0.7321428571428571
$ws.Range("J12").Value = 0.1964285714285714
$ws.Range("K12").Value = 0.008928571428571428
$ws.Range("L12").Value = 0.02678571428571428
$ws.Range("S12").Value = 0.03571428571428571
# Row 13
$ws.Range("G13").Value = 0.7727272727272727
$ws.Range("J13").Value = 0.1363636363636364
# Row 15
$ws.Range("F15").Value = 0.0196078431372549
$ws.Range("H15").Value = 0.142156862745098
$ws.Range("I15").Value = 0.08333333333333333
$ws.Range("J15").Value = 0.4117647058823529
$ws.Range("K15").Value = 0.009803921568627451
$ws.Range("M15").Value = 0.01470588235294118
$ws.Range("O15").Value = 0.04901960784313725
$ws.Range("S15").Value = 0.2696078431372549
# Row 16
$ws.Range("F16").Value = 0.01530612244897959
$ws.Range("H16").Value = 0.1530612244897959
$ws.Range("I16").Value = 0.06122448979591837
$ws.Range("J16").Value = 0.4234693877551021
$ws.Range("K16").Value = 0.08673469387755102
$ws.Range("M16").Value = 0.04081632653061224
$ws.Range("O16").Value = 0.05102040816326531
$ws.Range("S16").Value = 0.1683673469387755
# Row 17
$ws.Range("F17").Value = 0.018140589569161
$ws.Range("H17").Value = 0.1972789115646258
$ws.Range("I17").Value = 0.1383219954648526
$ws.Range("J17").Value = 0.4172335600907029
$ws.Range("K17").Value = 0.07256235827664399
$ws.Range("M17").Value = 0.02040816326530612
$ws.Range("N17").Value = 0.002267573696145125
$ws.Range("O17").Value = 0.06122448979591837
$ws.Range("S17").Value = 0.07256235827664399
# Row 18
$ws.Range("F18").Value = 0.01507537688442211
$ws.Range("H18").Value = 0.1758793969849246
$ws.Range("I18").Value = 0.1055276381909548
$ws.Range("J18").Value = 0.4321608040201005
$ws.Range("K18").Value = 0.06532663316582915
$ws.Range("M18").Value = 0.01507537688442211
$ws.Range("O18").Value = 0.06532663316582915
$ws.Range("S18").Value = 0.1256281407035176
# Row 19
$ws.Range("F19").Value = 0.01724137931034483
$ws.Range("H19").Value = 0.2270114942528736
$ws.Range("I19").Value = 0.09386973180076628
$ws.Range("J19").Value = 0.3927203065134099
$ws.Range("K19").Value = 0.08908045977011494
$ws.Range("M19").Value = 0.0210727969348659
$ws.Range("O19").Value = 0.07183908045977011
$ws.Range("S19").Value = 0.08716475095785441

Write-Host "Updated 106 cell values in the transition matrix"
